$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing used range to drop stale formatting/content before rewriting
$ws.Range("A1:K20").Clear()

$ws.Range("A1").Value = 'code'
$ws.Range("B1").Value = 'name'
$ws.Range("C1").Value = 'family'
$ws.Range("D1").Value = 'native_name'
$ws.Range("E1").Value = 'is_active'
$ws.Range("F1").Value = 'cr_by'
$ws.Range("G1").Value = 'cr_dtimes'
$ws.Range("H1").Value = 'upd_by'
$ws.Range("I1").Value = 'upd_dtimes'
$ws.Range("J1").Value = 'is_deleted'
$ws.Range("K1").Value = 'del_dtimes'
$ws.Range("A2").Value = 'eng'
$ws.Range("B2").Value = 'English'
$ws.Range("C2").Value = 'Indo-European'
$ws.Range("D2").Value = 'English'
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 'admin'
$ws.Range("G2").Value = 45589.5135756683
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = 'NULL'
$ws.Range("I2").Value = 'NULL'
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = 'NULL'
$ws.Range("A3").Value = 'fra'
$ws.Range("B3").Value = 'French'
$ws.Range("C3").Value = 'Indo-European'
$ws.Range("D3").Value = 'français'
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 'admin'
$ws.Range("G3").Value = 45589.5135756683
$ws.Range("G3").NumberFormat = "mm:ss.0"
$ws.Range("H3").Value = 'NULL'
$ws.Range("I3").Value = 'NULL'
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = 'NULL'
$ws.Range("A4").Value = 'ara'
$ws.Range("B4").Value = 'Arabic'
$ws.Range("C4").Value = 'الهندو أوروبية'
$ws.Range("D4").Value = 'Arabic'
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 'admin'
$ws.Range("G4").Value = 45589.5135756683
$ws.Range("G4").NumberFormat = "mm:ss.0"
$ws.Range("H4").Value = 'NULL'
$ws.Range("I4").Value = 'NULL'
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = 'NULL'
$ws.Range("A5").Value = 'kan'
$ws.Range("B5").Value = 'ಕನ್ನಡ'
$ws.Range("C5").Value = 'ಇಂಡೋ-ಯುರೋಪಿಯನ್'
$ws.Range("D5").Value = 'Kannada'
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 'admin'
$ws.Range("G5").Value = 45589.5135756683
$ws.Range("G5").NumberFormat = "mm:ss.0"
$ws.Range("H5").Value = 'NULL'
$ws.Range("I5").Value = 'NULL'
$ws.Range("J5").Value = $false
$ws.Range("K5").Value = 'NULL'
$ws.Range("A6").Value = 'hin'
$ws.Range("B6").Value = 'हिन्दी'
$ws.Range("C6").Value = 'भारोपीय'
$ws.Range("D6").Value = 'Hindi'
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 'admin'
$ws.Range("G6").Value = 45589.5135756683
$ws.Range("G6").NumberFormat = "mm:ss.0"
$ws.Range("H6").Value = 'NULL'
$ws.Range("I6").Value = 'NULL'
$ws.Range("J6").Value = $false
$ws.Range("K6").Value = 'NULL'
$ws.Range("A7").Value = 'tam'
$ws.Range("B7").Value = 'தமிழ்'
$ws.Range("C7").Value = 'இந்தோ-ஐரோப்பிய'
$ws.Range("D7").Value = 'Tamil'
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = 'admin'
$ws.Range("G7").Value = 45589.5135756683
$ws.Range("G7").NumberFormat = "mm:ss.0"
$ws.Range("H7").Value = 'NULL'
$ws.Range("I7").Value = 'NULL'
$ws.Range("J7").Value = $false
$ws.Range("K7").Value = 'NULL'
$ws.Range("A8").Value = 'es'
$ws.Range("B8").Value = 'Spanish'
$ws.Range("C8").Value = 'Indo-European'
$ws.Range("D8").Value = 'Spanish'
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = 'admin'
$ws.Range("G8").Value = 45589.5135756683
$ws.Range("G8").NumberFormat = "mm:ss.0"
$ws.Range("H8").Value = 'NULL'
$ws.Range("I8").Value = 'NULL'
$ws.Range("J8").Value = $false
$ws.Range("K8").Value = 'NULL'

$ws.Range("A1").Select() | Out-Null